$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix player name loop bug: row 4 had wrong player name
$ws.Range("A4").Value = "ranjan ramanayake"

# Add new "PLAYER NUMBER" column header and a sample value
$ws.Range("G1").Value = "PLAYER NUMBER"
$ws.Range("G2").Value = "aa"
